$p = $ppt.ActivePresentation

# The deck had 11 slides; the last two ("Add a Slide Title - 4" and
# "Add a Slide Title - 5") are no longer wanted, bringing the total down
# to 9 slides.
for ($i = 0; $i -lt 2; $i++) {
    $last = $p.Slides.Count
    $p.Slides.Item($last).Delete()
}
